# Insert a new weekly price record at row 9 (Fruta / hortaliza, semanal).
# This shifts the existing rows 9-75 down to 10-76, matching the
# "Fruta, Terminal Hortofrutícola Agro Chillán - Ciruela" subconjunto diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 9 (and everything below it) down by one row.
$ws.Rows.Item(9).Insert()

# Populate the newly-opened row 9 with the new weekly record.
$ws.Cells.Item(9, 1).Value  = 7
$ws.Cells.Item(9, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value  = "Ñuble"
$ws.Cells.Item(9, 4).Value  = 44959
$ws.Cells.Item(9, 5).Value  = 16
$ws.Cells.Item(9, 6).Value  = "Fruta"
$ws.Cells.Item(9, 7).Value  = 100103
$ws.Cells.Item(9, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(9, 9).Value  = 100103002
$ws.Cells.Item(9, 10).Value = "Ciruela"
$ws.Cells.Item(9, 11).Value = "Black Amber"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 60
$ws.Cells.Item(9, 14).Value = 10000
$ws.Cells.Item(9, 15).Value = 11000
$ws.Cells.Item(9, 16).Value = 10500
$ws.Cells.Item(9, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(9, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(9, 19).Value = 583
$ws.Cells.Item(9, 20).Value = 18
